$d = $word.ActiveDocument
# delete existing _GoBack
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$ftr = $d.Sections(1).Footers.Item(1)
$find = $ftr.Range.Find
$find.ClearFormatting()
$found = $find.Execute("12/5/2018", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $find.Parent.Start
$dateEnd = $find.Parent.End

$fiveRange = $ftr.Range.Duplicate
$fiveRange.SetRange($dateStart + 3, $dateStart + 4)
$fiveRange.Text = "9"

# now add a new _GoBack bookmark right after the "9" (collapsed range)
$bkRange = $ftr.Range.Duplicate
$bkRange.SetRange($dateStart + 4, $dateStart + 4)
$d.Bookmarks.Add("_GoBack", $bkRange)
Write-Host "GoBack exists:" $d.Bookmarks.Exists("_GoBack")
